$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J5").Value = 3.6
$ws.Range("M5").Value = 1.11
$ws.Range("N5").Value = 6.5
$ws.Range("O5").Value = 1.53
$ws.Range("P5").Value = 2.38
$ws.Range("Q5").Value = 2.7
$ws.Range("R5").Value = 1.44
$ws.Range("U5").Value = 2.1
$ws.Range("V5").Value = 1.67
$ws.Range("AF5").Value = 67
$ws.Range("AI5").Value = 12
$ws.Range("AN5").Value = 4.5
$ws.Range("AV5").Value = 67
$ws.Range("M6").Value = 1.17
$ws.Range("N6").Value = 5
$ws.Range("AC6").Value = 5
$ws.Range("AI6").Value = 13
$ws.Range("AS6").Value = 451
$ws.Range("AU6").Value = 10
$ws.Range("AW6").Value = 4.75
$ws.Range("O7").Value = 1.53
$ws.Range("P7").Value = 2.38
$ws.Range("Q7").Value = 2.7
$ws.Range("R7").Value = 1.44
$ws.Range("AH7").Value = 9
$ws.Range("M8").Value = 1.17
$ws.Range("N8").Value = 5
$ws.Range("W8").Value = 6
$ws.Range("AA8").Value = 29
$ws.Range("AC8").Value = 5
$ws.Range("K9").Value = 1.95
$ws.Range("N9").Value = 7.5
$ws.Range("U9").Value = 2.1
$ws.Range("V9").Value = 1.67
$ws.Range("X9").Value = 9
$ws.Range("AA9").Value = 21
$ws.Range("AB9").Value = 41
$ws.Range("AR9").Value = 81
$ws.Range("AS9").Value = 251
$ws.Range("O10").Value = 1.3
$ws.Range("P10").Value = 3.5
$ws.Range("Q10").Value = 2.03
$ws.Range("R10").Value = 1.83
$ws.Range("G13").Value = 2.65
$ws.Range("H13").Value = 3.4
$ws.Range("J13").Value = 3.15
$ws.Range("K13").Value = 2.25
$ws.Range("L13").Value = 2.9
$ws.Range("N13").Value = 9
$ws.Range("S13").Value = 1.3
$ws.Range("T13").Value = 3.25
$ws.Range("V13").Value = 2.67
$ws.Range("W13").Value = 13.5
$ws.Range("X13").Value = 17.5
$ws.Range("Z13").Value = 35
$ws.Range("AA13").Value = 19
$ws.Range("AC13").Value = 9
$ws.Range("AD13").Value = 7.3
$ws.Range("AE13").Value = 10.25
$ws.Range("AF13").Value = 30
$ws.Range("AK13").Value = 30
$ws.Range("AM13").Value = 18.5
$ws.Range("AN13").Value = 5.1
$ws.Range("AO13").Value = 14
$ws.Range("AP13").Value = 17
$ws.Range("AQ13").Value = 55
$ws.Range("AT13").Value = 3.25
$ws.Range("AU13").Value = 6
$ws.Range("AV13").Value = 37
$ws.Range("AY13").Value = 15.5
$ws.Range("BC13").Value = 450
$ws.Range("G27").Value = 2.63
$ws.Range("I27").Value = 2.75
$ws.Range("M27").Value = 1.1
$ws.Range("N27").Value = 7
$ws.Range("X27").Value = 11
$ws.Range("Z27").Value = 26
$ws.Range("AL27").Value = 29
$ws.Range("AW27").Value = 4.75
$ws.Range("O28").Value = 1.29
$ws.Range("P28").Value = 3.5
$ws.Range("Q28").Value = 2
$ws.Range("R28").Value = 1.85
$ws.Range("O36").Value = 1.14
$ws.Range("P36").Value = 5.5
$ws.Range("O42").Value = 1.3
$ws.Range("P42").Value = 3.5
$ws.Range("Q42").Value = 1.98
$ws.Range("R42").Value = 1.88
$ws.Range("G47").Value = 3.5
$ws.Range("I47").Value = 1.91
$ws.Range("K47").Value = 2.1
$ws.Range("L47").Value = 2.63
$ws.Range("M47").Value = 1.06
$ws.Range("N47").Value = 10
$ws.Range("Q47").Value = 2.03
$ws.Range("R47").Value = 1.83
$ws.Range("U47").Value = 1.83
$ws.Range("V47").Value = 1.83
$ws.Range("AA47").Value = 29
$ws.Range("AD47").Value = 6.5
$ws.Range("AH47").Value = 7
$ws.Range("AL47").Value = 17
$ws.Range("AM47").Value = 29
$ws.Range("AX47").Value = 11
$ws.Range("AZ47").Value = 41
$ws.Range("O51").Value = 1.22
$ws.Range("P51").Value = 4
$ws.Range("Q51").Value = 1.73
$ws.Range("R51").Value = 2.08
$ws.Range("Y54").Value = 6
$ws.Range("G56").Value = 2.05
$ws.Range("H56").Value = 3.2
$ws.Range("I56").Value = 4.1
$ws.Range("U56").Value = 1.95
$ws.Range("V56").Value = 1.8
$ws.Range("W56").Value = 6.5
$ws.Range("X56").Value = 9
$ws.Range("Z56").Value = 17
$ws.Range("AA56").Value = 17
$ws.Range("AO56").Value = 11
$ws.Range("Q58").Value = 1.67
$ws.Range("R58").Value = 2.15
$ws.Range("I59").Value = 3.2
$ws.Range("J59").Value = 2.6
$ws.Range("L59").Value = 3.4
$ws.Range("N59").Value = 19
$ws.Range("S59").Value = 1.22
$ws.Range("T59").Value = 4
$ws.Range("X59").Value = 15
$ws.Range("Y59").Value = 9.5
$ws.Range("AB59").Value = 17
$ws.Range("AN59").Value = 4.75
$ws.Range("AO59").Value = 11
$ws.Range("AQ59").Value = 34
$ws.Range("AT59").Value = 4
$ws.Range("AU59").Value = 6.5
$ws.Range("AZ59").Value = 41
$ws.Range("BB59").Value = 81
$ws.Range("G60").Value = 4.85
$ws.Range("W60").Value = 13
$ws.Range("AH60").Value = 7.3
$ws.Range("AN60").Value = 6.8
$ws.Range("AW60").Value = 3.55
$ws.Range("G61").Value = 2.25
$ws.Range("I61").Value = 3.5
$ws.Range("J61").Value = 2.88
$ws.Range("K61").Value = 2.1
$ws.Range("N61").Value = 8.5
$ws.Range("X61").Value = 10
$ws.Range("AC61").Value = 8.5
$ws.Range("AD61").Value = 6
$ws.Range("AH61").Value = 10
$ws.Range("AJ61").Value = 13
$ws.Range("AK61").Value = 41
$ws.Range("AO61").Value = 12
$ws.Range("AS61").Value = 151
$ws.Range("H63").Value = 4.05
$ws.Range("I63").Value = 5
$ws.Range("J63").Value = 2.05
$ws.Range("K63").Value = 2.32
$ws.Range("O63").Value = 1.2
$ws.Range("P63").Value = 3.6
$ws.Range("Q63").Value = 1.62
$ws.Range("R63").Value = 2.02
$ws.Range("S63").Value = 1.31
$ws.Range("T63").Value = 3.26
$ws.Range("U63").Value = 1.7
$ws.Range("V63").Value = 1.91
$ws.Range("W63").Value = 7.8
$ws.Range("X63").Value = 7.9
$ws.Range("AA63").Value = 11.75
$ws.Range("AC63").Value = 13
$ws.Range("AD63").Value = 8
$ws.Range("AJ63").Value = 16
$ws.Range("AL63").Value = 45
$ws.Range("AO63").Value = 7.2
$ws.Range("AR63").Value = 45
$ws.Range("AS63").Value = 175
$ws.Range("AT63").Value = 3.05
$ws.Range("U67").Value = 1.57
$ws.Range("U68").Value = 1.57
$ws.Range("U69").Value = 1.57
$ws.Range("U70").Value = 1.5
$ws.Range("G78").Value = 1.62
$ws.Range("H78").Value = 3.6
$ws.Range("I78").Value = 5.5
$ws.Range("J78").Value = 2.25
$ws.Range("K78").Value = 2.2
$ws.Range("L78").Value = 6
$ws.Range("M78").Value = 1.06
$ws.Range("N78").Value = 10
$ws.Range("O78").Value = 1.33
$ws.Range("P78").Value = 3.25
$ws.Range("Q78").Value = 2.08
$ws.Range("R78").Value = 1.73
$ws.Range("U78").Value = 2.1
$ws.Range("V78").Value = 1.67
$ws.Range("AH78").Value = 13
$ws.Range("AJ78").Value = 19
$ws.Range("AO78").Value = 8.5
$ws.Range("AW78").Value = 7
$ws.Range("AZ78").Value = 126
$ws.Range("BB78").Value = 351
